# Apply the "OCM audio video play back report" query documentation edit
# to the Queries sheet: add two long SQL text cells (F2, G2) with wrapped,
# styled formatting, widen their columns, grow row 2 to fit, and update the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Queries")
$ws.Activate()

$sql1 = @"
SELECT M.AgentID as [Agent ID],A.AgentName as[Agent Name],A.TeamName as [Team Name],A.SupervisorName as[Supervisor Name],sum(ACDCalls) AS [Total Interaction],[dbo].[SECONDSTOhhmmss](SUM([TotalInteractionTime])) as [Total Interaction Time],
[dbo].[SECONDSTOhhmmss](SUM([TotalInteractionTime])/nullif(SUM(ACDCalls),0))AS [Avg Interaction Time],SUM([TotalChat]) as [Total Chat],
[dbo].[SECONDSTOhhmmss](SUM([TotalChatTime])) as[Total Chat Time],[dbo].[SECONDSTOhhmmss](SUM([TotalChatTime])/nullif(sum([TotalChat]),0))AS [Avg Chat Time],SUM([TotalAudioIP]) as [Total Audio IP],[dbo].[SECONDSTOhhmmss](sum([TotalAudioIPTime])) [Total AudioIP Time],
[dbo].[SECONDSTOhhmmss](SUM([TotalAudioIPTime])/nullif(SUM([TotalAudioIP]),0))AS [Avg AudioIP Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalInteractionTime)/nullif(SUM(ACDCalls),0))AS [Avg Talk Time],[dbo].[SECONDSTOhhmmss](SUM(TotalAfterCallTime)) as [Total After Call Time],[dbo].[SECONDSTOhhmmss](SUM(TotalAvailTime)) as [Total Avail Time],[dbo].[SECONDSTOhhmmss](SUM(TotalAuxTime)) as [Total AUX Time],SUM(ExtensionCalls) AS [Extension Interaction],[dbo].[SECONDSTOhhmmss](sum([TotalExtensionTime])) as [Total Extension Time],
[dbo].[SECONDSTOhhmmss](sum(TotalExtensionTime)/nullif(sum(ExtensionCalls),0))AS [Avg Extension Time],[dbo].[SECONDSTOhhmmss](SUM(TotalStaffedTime)) as [Total Time Staffed],[dbo].[SECONDSTOhhmmss](SUM(TotalHoldTime)) as [Total Hold Time] 
FROM [OCM_AgentHistoricalReport] M WITH(NOLOCK)
INNER JOIN fn_AgentHierarchy('na','1','1') A ON  A.[AgentId]=M.[AgentID]
WHERE [ReportDateTime]>='ReportBeforeDate' and [ReportDateTime]<='ReportAfterDate'
GROUP BY M.[AgentID], A.[AgentName],A.[TeamName],A.[SupervisorName]
Order by [Agent Name];
"@

$sql2 = @"
SELECT Dateint AS [Date],AgentID as [Agent ID],ISNULL(A.FirstName,'')+' '+ ISNULL(A.LastName,'') AS [Agent Name],ISNULL(C.TeamName,' ') AS TeamName,
ISNULL(B.FirstName,'NA')+' '+ ISNULL(B.LastName,'') AS SupervisorName,SUM(ACDCalls) AS TotalInteraction,
[dbo].[SECONDSTOhhmmss](SUM([TotalInteractionTime])) [Total Interaction Time],
[dbo].[SECONDSTOhhmmss](SUM([TotalInteractionTime])/nullif(SUM(ACDCalls),0))AS [Avg Interaction Time],
SUM([TotalChat]) as [Total Chat],[dbo].[SECONDSTOhhmmss](SUM([TotalChatTime])) as [Total Chat Time],
[dbo].[SECONDSTOhhmmss](SUM([TotalChatTime])/nullif(sum([TotalChat]),0))AS [Avg Chat Time],
SUM([TotalAudioIP]) [Total Audio IP],[dbo].[SECONDSTOhhmmss](sum([TotalAudioIPTime])) [Total AudioIP Time],
[dbo].[SECONDSTOhhmmss](SUM([TotalAudioIPTime])/nullif(SUM([TotalAudioIP]),0))AS [Avg AudioIP Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalInteractionTime)/nullif(SUM(ACDCalls),0))AS [Avg Talk Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalAfterCallTime)) as [Total After Call Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalAvailTime)) as [Total Avail Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalAuxTime)) as [Total Aux Time],
SUM(ExtensionCalls) AS [Extension Interaction],
[dbo].[SECONDSTOhhmmss](sum([TotalExtensionTime])) [Total Extension Time],
[dbo].[SECONDSTOhhmmss](sum(TotalExtensionTime)/nullif(sum(ExtensionCalls),0))AS [Avg Extension Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalStaffedTime)) as [Total Time Staffed],
[dbo].[SECONDSTOhhmmss](SUM(TotalHoldTime)) as [Total Hold Time]  
FROM [OCM_AgentHistoricalReport] M WITH(NOLOCK)
LEFT JOIN [AGT_Agent] A WITH(NOLOCK)  ON A.AvayaLoginID = M.[AgentID] 
LEFT JOIN [AGT_Agent] B WITH(NOLOCK)  ON A.[PrimarySupervisorID]=B.ID LEFT JOIN [AGT_Teams] C WITH(NOLOCK) ON C.TeamID = A.TeamID
LEFT JOIN[dbo].[AGT_Teams] P WITH(NOLOCK) ON C.ParentID = P.TeamID 
WHERE [ReportDateTime]>='ReportBeforeDate' AND [ReportDateTime]<='ReportAfterDate' AND [AgentID] LIKE 'AgentIdCapturedFromUI'  
GROUP BY  [Dateint], [AgentID],B.[FirstName],B.[LastName],C.TeamName,A.FirstName,A.LastName ORDER BY [Dateint] ASC;
"@

# New shared-string cells holding the SQL query text used by the report.
$ws.Range("F2").Value = $sql1
$ws.Range("G2").Value = $sql2

# F2: wrap text only.
$ws.Range("F2").WrapText = $true

# G2: centered (horizontal + vertical) and wrapped.
$ws.Range("G2").WrapText = $true
$ws.Range("G2").HorizontalAlignment = -4108
$ws.Range("G2").VerticalAlignment = -4108

# Widen columns F and G to fit the new long text.
$ws.Columns.Item(6).ColumnWidth = 82
$ws.Columns.Item(7).ColumnWidth = 56.6666666667

# Grow row 2 to the maximum row height to accommodate the wrapped text.
$ws.Rows.Item(2).RowHeight = 409.5

# Update view: scroll/selection now centers on the newly added content.
$ws.Range("G2").Select() | Out-Null
